{"js": "// Append a new bulleted list item after the last paragraph of the document\n// (\"Frame choosing lais\u00b2 ...\"), continuing the same list (Listenabsatz\n// style, numId 1, ilvl 0), as described by the diff / commit message\n// (\"Working on difussion arrays\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The paragraph that currently ends the document body.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newText =\n  \"My community tracing might not take the following case into account: \" +\n  \"2 communities are overlapping. In the next time step both communities \" +\n  \"expand and a 3. community spawns. This 3. Community consists only of \" +\n  \"nodes that are members of one of the two previous communities. Investigate\";\n\n// Inserting \"After\" the existing list paragraph creates a new paragraph\n// that inherits its paragraph style / list formatting (Listenabsatz,\n// numId 1, ilvl 0), matching the target OOXML in the diff.\nlastParagraph.insertParagraph(newText, \"After\");\n\nawait context.sync();\n", "ps1": "# Append a new bulleted list item after the last paragraph of the document\n# (\"Frame choosing lais\u00b2 ...\"), continuing the same list (Listenabsatz\n# style, numId 1, ilvl 0), as described by the diff / commit message\n# (\"Working on difussion arrays\").\n\n$d = $word.ActiveDocument\n\n# The paragraph that currently ends the document body.\n$lastIndex = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs($lastIndex)\n\n$newText = \"My community tracing might not take the following case into account: 2 communities are overlapping. In the next time step both communities expand and a 3. community spawns. This 3. Community consists only of nodes that are members of one of the two previous communities. Investigate\"\n\n# Inserting a paragraph mark right after the existing list paragraph creates\n# a new paragraph that inherits its paragraph style / list formatting\n# (Listenabsatz, numId 1, ilvl 0), matching the target OOXML in the diff.\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$newParagraph.Range.Text = $newText\n"}
